$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new column F header (string value, matches new shared string "11_03_2024")
$ws.Range("F1").Value = "11_03_2024"

# Fill in the F column values for rows 2-6 (row 6 cell already existed with a style)
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 4

# Update the selected cell to F7, matching the saved selection in the diff
$ws.Range("F7").Select()

$wb.Save()
